# ------------------------------------------------------------------
# Adds the "2022-Q4" quarter to the workbook:
#   1. Inserts a new "2022-Q4" worksheet (with its fund-holdings table)
#      right after "总计" and before the existing "2022-Q3" sheet.
#   2. Inserts a matching summary row into the "总计" (totals) sheet.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---- 1. Create the new "2022-Q4" sheet, placed before "2022-Q3" ----
$totals = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Item("2022-Q3")
$ws = $wb.Worksheets.Add($q3)
$ws.Name = "2022-Q4"

# Match page margins used by the other quarter sheets (0.75/0.75/1/1/.5/.5 in)
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---- Fund-holdings data for 2022-Q4 ----
    $ws.Range("B1").Value = "基金代码"
    $ws.Range("C1").Value = "基金名称"
    $ws.Range("D1").Value = "基金规模"
    $ws.Range("E1").Value = "股票总仓位"
    $ws.Range("F1").Value = "仓位占比"
    $ws.Range("G1").Value = "持有市值(亿元)"
    $ws.Range("H1").Value = "仓位排名"
    $ws.Range("A2").Value = 0
    $ws.Range("B2").Value = "'501077"
    $ws.Range("C2").Value = "富国创新企业灵活配置混合（LOF）A"
    $ws.Range("D2").Value = "'8.06"
    $ws.Range("E2").Value = "'88.91"
    $ws.Range("F2").Value = "'2.71"
    $ws.Range("G2").Value = "'0.2184"
    $ws.Range("H2").Value = 9
    $ws.Range("A3").Value = 1
    $ws.Range("B3").Value = "'002450"
    $ws.Range("C3").Value = "平安睿享文娱灵活配置混合A"
    $ws.Range("D3").Value = "'3.85"
    $ws.Range("E3").Value = "'94.03"
    $ws.Range("F3").Value = "'4.03"
    $ws.Range("G3").Value = "'0.1552"
    $ws.Range("H3").Value = 7
    $ws.Range("A4").Value = 2
    $ws.Range("B4").Value = "'006101"
    $ws.Range("C4").Value = "平安优势产业灵活配置混合C"
    $ws.Range("D4").Value = "'2.85"
    $ws.Range("E4").Value = "'94.90"
    $ws.Range("F4").Value = "'3.51"
    $ws.Range("G4").Value = "'0.1000"
    $ws.Range("H4").Value = 8
    $ws.Range("A5").Value = 3
    $ws.Range("B5").Value = "'010126"
    $ws.Range("C5").Value = "平安价值成长混合A"
    $ws.Range("D5").Value = "'3.21"
    $ws.Range("E5").Value = "'94.03"
    $ws.Range("F5").Value = "'3.09"
    $ws.Range("G5").Value = "'0.0992"
    $ws.Range("H5").Value = 9
    $ws.Range("A6").Value = 4
    $ws.Range("B6").Value = "'011828"
    $ws.Range("C6").Value = "平安睿享成长混合A"
    $ws.Range("D6").Value = "'2.42"
    $ws.Range("E6").Value = "'92.96"
    $ws.Range("F6").Value = "'3.30"
    $ws.Range("G6").Value = "'0.0799"
    $ws.Range("H6").Value = 9
    $ws.Range("A7").Value = 5
    $ws.Range("B7").Value = "'002451"
    $ws.Range("C7").Value = "平安睿享文娱灵活配置混合C"
    $ws.Range("D7").Value = "'1.97"
    $ws.Range("E7").Value = "'94.03"
    $ws.Range("F7").Value = "'4.03"
    $ws.Range("G7").Value = "'0.0794"
    $ws.Range("H7").Value = 7
    $ws.Range("A8").Value = 6
    $ws.Range("B8").Value = "'006100"
    $ws.Range("C8").Value = "平安优势产业灵活配置混合A"
    $ws.Range("D8").Value = "'2.15"
    $ws.Range("E8").Value = "'94.90"
    $ws.Range("F8").Value = "'3.51"
    $ws.Range("G8").Value = "'0.0755"
    $ws.Range("H8").Value = 8
    $ws.Range("A9").Value = 7
    $ws.Range("B9").Value = "'013687"
    $ws.Range("C9").Value = "平安成长龙头1年持有混合A"
    $ws.Range("D9").Value = "'1.21"
    $ws.Range("E9").Value = "'94.87"
    $ws.Range("F9").Value = "'4.02"
    $ws.Range("G9").Value = "'0.0486"
    $ws.Range("H9").Value = 8
    $ws.Range("A10").Value = 8
    $ws.Range("B10").Value = "'010127"
    $ws.Range("C10").Value = "平安价值成长混合C"
    $ws.Range("D10").Value = "'1.33"
    $ws.Range("E10").Value = "'94.03"
    $ws.Range("F10").Value = "'3.09"
    $ws.Range("G10").Value = "'0.0411"
    $ws.Range("H10").Value = 9
    $ws.Range("A11").Value = 9
    $ws.Range("B11").Value = "'011829"
    $ws.Range("C11").Value = "平安睿享成长混合C"
    $ws.Range("D11").Value = "'1.09"
    $ws.Range("E11").Value = "'92.96"
    $ws.Range("F11").Value = "'3.30"
    $ws.Range("G11").Value = "'0.0360"
    $ws.Range("H11").Value = 9
    $ws.Range("A12").Value = 10
    $ws.Range("B12").Value = "'010642"
    $ws.Range("C12").Value = "农银汇理瑞祥一年持有期混合"
    $ws.Range("D12").Value = "'2.48"
    $ws.Range("E12").Value = "'21.31"
    $ws.Range("F12").Value = "'1.40"
    $ws.Range("G12").Value = "'0.0347"
    $ws.Range("H12").Value = 7
    $ws.Range("A13").Value = 11
    $ws.Range("B13").Value = "'200001"
    $ws.Range("C13").Value = "长城久恒灵活配置混合"
    $ws.Range("D13").Value = "'0.85"
    $ws.Range("E13").Value = "'94.35"
    $ws.Range("F13").Value = "'2.59"
    $ws.Range("G13").Value = "'0.0220"
    $ws.Range("H13").Value = 10
    $ws.Range("A14").Value = 12
    $ws.Range("B14").Value = "'013688"
    $ws.Range("C14").Value = "平安成长龙头1年持有混合C"
    $ws.Range("D14").Value = "'0.51"
    $ws.Range("E14").Value = "'94.87"
    $ws.Range("F14").Value = "'4.02"
    $ws.Range("G14").Value = "'0.0205"
    $ws.Range("H14").Value = 8
    $ws.Range("A15").Value = 13
    $ws.Range("B15").Value = "'007894"
    $ws.Range("C15").Value = "平安估值精选混合C"
    $ws.Range("D15").Value = "'0.50"
    $ws.Range("E15").Value = "'94.39"
    $ws.Range("F15").Value = "'3.54"
    $ws.Range("G15").Value = "'0.0177"
    $ws.Range("H15").Value = 8
    $ws.Range("A16").Value = 14
    $ws.Range("B16").Value = "'007893"
    $ws.Range("C16").Value = "平安估值精选混合A"
    $ws.Range("D16").Value = "'0.34"
    $ws.Range("E16").Value = "'94.39"
    $ws.Range("F16").Value = "'3.54"
    $ws.Range("G16").Value = "'0.0120"
    $ws.Range("H16").Value = 8
    $ws.Range("A17").Value = 15
    $ws.Range("B17").Value = "'006721"
    $ws.Range("C17").Value = "平安核心优势混合C"
    $ws.Range("D17").Value = "'0.17"
    $ws.Range("E17").Value = "'89.85"
    $ws.Range("F17").Value = "'4.44"
    $ws.Range("G17").Value = "'0.0075"
    $ws.Range("H17").Value = 7
    $ws.Range("A18").Value = 16
    $ws.Range("B18").Value = "'006720"
    $ws.Range("C18").Value = "平安核心优势混合A"
    $ws.Range("D18").Value = "'0.06"
    $ws.Range("E18").Value = "'89.85"
    $ws.Range("F18").Value = "'4.44"
    $ws.Range("G18").Value = "'0.0027"
    $ws.Range("H18").Value = 7
    $ws.Range("A19").Value = 17
    $ws.Range("B19").Value = "'015849"
    $ws.Range("C19").Value = "富国创新企业灵活配置混合（LOF）C"
    $ws.Range("D19").Value = "'0.01"
    $ws.Range("E19").Value = "'88.91"
    $ws.Range("F19").Value = "'2.71"
    $ws.Range("G19").Value = "'0.0003"
    $ws.Range("H19").Value = 9

# Header row (B1:H1) + index column (A2:A19) use the bold/centered/bordered
# style already used for headers elsewhere in the workbook.
$headerStyle = $totals.Range("B1").Style
$ws.Range("B1:H1").Style = $headerStyle
$ws.Range("A2:A19").Style = $headerStyle

# ---- 2. Insert the 2022-Q4 summary row into the "总计" sheet ----
# Push the existing rows (2022-Q3, 2022-Q2) down one row first, writing
# explicit values instead of Rows.Insert so the already-correct cached
# numbers aren't disturbed.
$indexStyle = $totals.Range("A2").Style

$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q2"
$totals.Range("C4").Value = 3
$totals.Range("D4").Value = 0.28
$totals.Range("A4").Style = $indexStyle

$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q3"
$totals.Range("C3").Value = 17
$totals.Range("D3").Value = 1.42
$totals.Range("A3").Style = $indexStyle

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 18
$totals.Range("D2").Value = 1.05
$totals.Range("A2").Style = $indexStyle
